$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.37372175046311
$ws.Range("C2").Value = 11.451355774266
$ws.Range("D2").Value = 5.965179928147954
$ws.Range("E2").Value = 16.54913201810116
$ws.Range("G2").Value = 3.601708615739476
$ws.Range("I2").Value = 17.370830946993
$ws.Range("O2").Value = 19.28612686431046

$ws.Range("B3").Value = 14.57208628456858
$ws.Range("C3").Value = 10.72204778027259
$ws.Range("D3").Value = 5.843154896066161
$ws.Range("E3").Value = 15.60490516367876
$ws.Range("G3").Value = 3.604515383009726
$ws.Range("I3").Value = 17.52909789762156
$ws.Range("O3").Value = 19.33277193092449

$ws.Range("B4").Value = 14.05711574745184
$ws.Range("C4").Value = 10.24651590419223
$ws.Range("D4").Value = 5.768801805894647
$ws.Range("E4").Value = 15.00029191135375
$ws.Range("G4").Value = 3.606326265246633
$ws.Range("I4").Value = 17.63276133484252
$ws.Range("O4").Value = 19.36991341533325

$ws.Range("B5").Value = 13.84174361239727
$ws.Range("C5").Value = 10.04575705902912
$ws.Range("D5").Value = 5.738692651106138
$ws.Range("E5").Value = 14.74792590966965
$ws.Range("G5").Value = 3.60708630347566
$ws.Range("I5").Value = 17.67662823464669
$ws.Range("O5").Value = 19.38716845643157

$ws.Range("B6").Value = 13.80565478508186
$ws.Range("C6").Value = 10.01199992413073
$ws.Range("D6").Value = 5.733705921085858
$ws.Range("E6").Value = 14.70566837845336
$ws.Range("G6").Value = 3.607213843787041
$ws.Range("I6").Value = 17.68401011073161
$ws.Range("O6").Value = 19.39016110448578

$ws.Range("B7").Value = 14.05423320259608
$ws.Range("C7").Value = 10.2438366439628
$ws.Range("D7").Value = 5.768394910447803
$ws.Range("E7").Value = 14.99691223866804
$ws.Range("G7").Value = 3.6063364258429
$ws.Range("I7").Value = 17.63334637731836
$ws.Range("O7").Value = 19.37013756582716

$ws.Range("B8").Value = 15.10215483671753
$ws.Range("C8").Value = 11.2056440255637
$ws.Range("D8").Value = 5.923015140349388
$ws.Range("E8").Value = 16.22886648137328
$ws.Range("G8").Value = 3.602658274965584
$ws.Range("I8").Value = 17.42405015354992
$ws.Range("O8").Value = 19.30043517717806

$ws.Range("B9").Value = 16.96920428669177
$ws.Range("C9").Value = 12.87199601542455
$ws.Range("D9").Value = 6.228734187518707
$ws.Range("E9").Value = 18.54299471116171
$ws.Range("G9").Value = 3.596136058385156
$ws.Range("I9").Value = 17.06546133971194
$ws.Range("O9").Value = 19.23194370882078

$ws.Range("B10").Value = 18.21864032852567
$ws.Range("C10").Value = 13.96296223996974
$ws.Range("D10").Value = 6.452255952971892
$ws.Range("E10").Value = 20.19849613350473
$ws.Range("G10").Value = 3.591759914635575
$ws.Range("I10").Value = 16.83413765173825
$ws.Range("O10").Value = 19.22410350419671

$ws.Range("B11").Value = 18.75933891070761
$ws.Range("C11").Value = 14.43051364001649
$ws.Range("D11").Value = 6.553183709969606
$ws.Range("E11").Value = 20.90950687465967
$ws.Range("G11").Value = 3.589858232041752
$ws.Range("I11").Value = 16.73600162333423
$ws.Range("O11").Value = 19.22991964844399

$ws.Range("B12").Value = 18.96003614463027
$ws.Range("C12").Value = 14.60344747943478
$ws.Range("D12").Value = 6.591253472331688
$ws.Range("E12").Value = 21.17273462188574
$ws.Range("G12").Value = 3.589150832337582
$ws.Range("I12").Value = 16.69987116074249
$ws.Range("O12").Value = 19.23348144102625

$ws.Range("B13").Value = 18.91699370217858
$ws.Range("C13").Value = 14.56638598059957
$ws.Range("D13").Value = 6.583061740411308
$ws.Range("E13").Value = 21.1163108477899
$ws.Range("G13").Value = 3.589302618767551
$ws.Range("I13").Value = 16.70760643879753
$ws.Range("O13").Value = 19.23265374749922

$ws.Range("B14").Value = 18.77593200122044
$ws.Range("C14").Value = 14.44482345065364
$ws.Range("D14").Value = 6.556318952933331
$ws.Range("E14").Value = 20.93128316030599
$ws.Range("G14").Value = 3.589799779239701
$ws.Range("I14").Value = 16.73300840833481
$ws.Range("O14").Value = 19.23018538489198

$ws.Range("B15").Value = 18.68899781214913
$ws.Range("C15").Value = 14.36982708454428
$ws.Range("D15").Value = 6.539917590526207
$ws.Range("E15").Value = 20.81716586305578
$ws.Range("G15").Value = 3.59010595938356
$ws.Range("I15").Value = 16.74870252851926
$ws.Range("O15").Value = 19.228850735688

$ws.Range("B16").Value = 18.18274082411529
$ws.Range("C16").Value = 13.93182963551521
$ws.Range("D16").Value = 6.445641363782595
$ws.Range("E16").Value = 20.15118606249291
$ws.Range("G16").Value = 3.591885979192407
$ws.Range("I16").Value = 16.84069478813083
$ws.Range("O16").Value = 19.22391311894114

$ws.Range("B17").Value = 17.86502246743362
$ws.Range("C17").Value = 13.65578253298937
$ws.Range("D17").Value = 6.387583042568138
$ws.Range("E17").Value = 19.73187512139006
$ws.Range("G17").Value = 3.59300071475193
$ws.Range("I17").Value = 16.89895401334582
$ws.Range("O17").Value = 19.2232949991598

$ws.Range("B18").Value = 17.67967757701557
$ws.Range("C18").Value = 13.4943019628114
$ws.Range("D18").Value = 6.354120829132503
$ws.Range("E18").Value = 19.48673521115121
$ws.Range("G18").Value = 3.593650267388636
$ws.Range("I18").Value = 16.93313021600751
$ws.Range("O18").Value = 19.22382213160332

$ws.Range("B19").Value = 17.61647864561873
$ws.Range("C19").Value = 13.43916231579524
$ws.Range("D19").Value = 6.342780604929844
$ws.Range("E19").Value = 19.4030532653738
$ws.Range("G19").Value = 3.593871637492479
$ws.Range("I19").Value = 16.94481591654738
$ws.Range("O19").Value = 19.2241518668106

$ws.Range("B20").Value = 17.89911397821292
$ws.Range("C20").Value = 13.68544813543526
$ws.Range("D20").Value = 6.393770854850403
$ws.Range("E20").Value = 19.77692142166106
$ws.Range("G20").Value = 3.592881181880105
$ws.Range("I20").Value = 16.89268310299178
$ws.Range("O20").Value = 19.22326936805781

$ws.Range("B21").Value = 18.81747576768983
$ws.Range("C21").Value = 14.48064094694426
$ws.Range("D21").Value = 6.564178327673541
$ws.Range("E21").Value = 20.98579328620168
$ws.Range("G21").Value = 3.589653406396741
$ws.Range("I21").Value = 16.72551914586926
$ws.Range("O21").Value = 19.23087343953663

$ws.Range("B22").Value = 19.39402693031205
$ws.Range("C22").Value = 14.97635152436459
$ws.Range("D22").Value = 6.674661360209797
$ws.Range("E22").Value = 21.74080783988571
$ws.Range("G22").Value = 3.587618010377652
$ws.Range("I22").Value = 16.62228595896577
$ws.Range("O22").Value = 19.24376975364773

$ws.Range("B23").Value = 19.08849536719137
$ws.Range("C23").Value = 14.71397190218229
$ws.Range("D23").Value = 6.615788693556396
$ws.Range("E23").Value = 21.34103814403468
$ws.Range("G23").Value = 3.58869758147146
$ws.Range("I23").Value = 16.67682891041219
$ws.Range("O23").Value = 19.2361586226724

$ws.Range("B24").Value = 17.88370956167406
$ws.Range("C24").Value = 13.67204496876018
$ws.Range("D24").Value = 6.390973602535889
$ws.Range("E24").Value = 19.75656867725348
$ws.Range("G24").Value = 3.592935195659513
$ws.Range("I24").Value = 16.89551605737939
$ws.Range("O24").Value = 19.22327820773487

$ws.Range("B25").Value = 16.48513354693592
$ws.Range("C25").Value = 12.4446703629532
$ws.Range("D25").Value = 6.146036951658545
$ws.Range("E25").Value = 17.89578715537306
$ws.Range("G25").Value = 3.597827102962322
$ws.Range("I25").Value = 17.15686628962754
$ws.Range("O25").Value = 19.24307161427025
